$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 414; existing rows 414-438 shift down to 415-439.
$ws.Rows.Item(414).Insert()

# Populate the newly inserted row 414 with the new weekly record.
$ws.Cells.Item(414, 1).Value = 4
$ws.Cells.Item(414, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(414, 3).Value = "Los Lagos"
$ws.Cells.Item(414, 4).Value = 44706
$ws.Cells.Item(414, 5).Value = 10
$ws.Cells.Item(414, 6).Value = "Fruta"
$ws.Cells.Item(414, 7).Value = 100102
$ws.Cells.Item(414, 8).Value = "Cítricos"
$ws.Cells.Item(414, 9).Value = 100102005
$ws.Cells.Item(414, 10).Value = "Naranja"
$ws.Cells.Item(414, 11).Value = "Valencia"
$ws.Cells.Item(414, 12).Value = "Segunda"
$ws.Cells.Item(414, 13).Value = 200
$ws.Cells.Item(414, 14).Value = 15000
$ws.Cells.Item(414, 15).Value = 15000
$ws.Cells.Item(414, 16).Value = 15000
$ws.Cells.Item(414, 17).Value = "$/caja 15 kilos empedrada"
$ws.Cells.Item(414, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(414, 19).Value = 1000
$ws.Cells.Item(414, 20).Value = 15
